# Actualización automática de tasas-transfi.xlsx

$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" note text ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.33 = 8807.08 pesos`n✅ 8807.08 pesos = 2.31 = 948.56 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- tasas: update the rate cells ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 429.995
$ws2.Range("O10").Value = 3787
$ws2.Range("N12").Value = 3815.99
$ws2.Range("O12").Value = 411
